# Rewrites the five placeholder "Chapter N" paragraphs into a repeating
# 4-paragraph cycle (origins / Ford / decades / modern), expanded to a
# total of 8 cycles (32 paragraphs), as produced by the RAG chain chatbot.

$d = $word.ActiveDocument

$origins = "The automobile's origins trace back to the late 19th century, when Karl Benz created the Benz Patent-Motorwagen in 1886. This invention marked the transition from horse-drawn carriages to self-propelled vehicles. Early automobiles were often expensive and unreliable, but they represented innovation in engineering and freedom of mobility."
$ford    = "Henry Ford's introduction of the moving assembly line in 1913 transformed production efficiency. The Ford Model T became a symbol of affordable transportation, making cars accessible to the average American. This innovation spread globally and shaped industrial production methods."
$decades = "Over decades, automotive technology evolved significantly. From the roaring 1920s with luxury classics to the muscle cars of the 1960s, and the oil crises of the 1970s, cars reflected cultural and economic shifts. Japanese manufacturers like Toyota and Honda rose in prominence due to reliability and fuel efficiency."
$modern  = "Modern cars are now safer and cleaner. Features like seatbelts, airbags, and electronic stability control are standard. Environmental regulations and emission standards continue to drive innovation toward cleaner energy and sustainability."

$cycle = @($origins, $ford, $decades, $modern)

$totalCycles = 8
$totalParagraphs = $totalCycles * 4

# Build the full ordered list of paragraph texts.
$all = @()
for ($i = 0; $i -lt $totalParagraphs; $i++) {
    $all += $cycle[$i % 4]
}

# The document currently has a title paragraph followed by 5 "Chapter N"
# placeholder paragraphs (index 2..6, 1-based). Rewrite those 5 in place
# (collapsing each one's multiple runs/line-breaks into a single run),
# then append the remaining paragraphs needed to reach the full cycle
# count.

$existingCount = 5

for ($i = 0; $i -lt $existingCount; $i++) {
    $p = $d.Paragraphs.Item($i + 2)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $all[$i]
}

$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastR = $lastP.Range

for ($i = $existingCount; $i -lt $all.Length; $i++) {
    $lastR.InsertParagraphAfter() | Out-Null
    $newP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newR = $newP.Range
    $newR.MoveEnd(1, -1) | Out-Null
    $newR.Text = $all[$i]
    $lastR = $newP.Range
}
